# Shift every existing booking date forward by exactly one year (365 days)
# and append a new Peak-Season / Limited row for the 2026 March 1 weekend.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 2-11: bump the Date column (A) forward by 365 days, keeping
# every other column untouched.
$ws.Range("A2").Value = 46015
$ws.Range("A3").Value = 46016
$ws.Range("A4").Value = 46022
$ws.Range("A5").Value = 46023
$ws.Range("A6").Value = 46024
$ws.Range("A7").Value = 46025
$ws.Range("A8").Value = 46063
$ws.Range("A9").Value = 46064
$ws.Range("A10").Value = 46069
$ws.Range("A11").Value = 46070

# New row 12 - mirrors row 7 (Limited / Peak Season) one year later.
# Copy A7's format (date-formatted, left-aligned) onto A12 before writing
# the value so the new cell reuses the existing date style.
$ws.Range("A7").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A12").Value = 46082

$ws.Range("B12").Value = "Limited"
$ws.Range("C12").Value = 18500
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = "Peak Season"

$ws.Range("J10").Select()
